$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (N) mirroring the per-row formatting already used
# by column M (2021), then apply a "0.0" number format to the numeric
# indicator rows so the new figures render with one decimal place.

# Row 2 (top border strip) - empty, formatting only
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

# Row 3 (year header)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2022

# Row 4
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 9.224468514531754
$ws.Range("N4").NumberFormat = "0.0"

# Row 5
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 4.6068543125097872
$ws.Range("N5").NumberFormat = "0.0"

# Row 6
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 13.543910285971602
$ws.Range("N6").NumberFormat = "0.0"

# Row 7
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 24.703327617190443
$ws.Range("N7").NumberFormat = "0.0"

# Row 8
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 28.608474183838851
$ws.Range("N8").NumberFormat = "0.0"

# Row 9
$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 20.904451081350146
$ws.Range("N9").NumberFormat = "0.0"

# Row 10
$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 26.720095429750884
$ws.Range("N10").NumberFormat = "0.0"

# Row 11
$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 27.704327204727914
$ws.Range("N11").NumberFormat = "0.0"

# Row 12
$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 25.731792255708452
$ws.Range("N12").NumberFormat = "0.0"

$ws.Range("Q5").Select()
